$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 58: "Se han asignado los permisos de forma restrictiva..." ---
# (set first so this string is appended to the shared-strings table before
# the row 57 one, matching the original authoring order)
$ws.Range("D58").Value = "SI"
$ws.Range("E53").Copy()
$ws.Range("E58").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E58").Value = "En system.sql y paquete de creación de usuarios"

# --- Row 57: "Se han asignado usuarios a los roles adecuadamente" ---
# Mark as answered "SI" and update the evidence note, matching the formatting
# used by the other "SI" evidence cells (italic black text instead of the
# italic blue "SEGURIDAD MARIO" placeholder note).
$ws.Range("D57").Value = "SI"
$ws.Range("E53").Copy()
$ws.Range("E57").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E57").Value = "En paquete de creación de usuarios"

$excel.CutCopyMode = 0

# --- Update the active selection / scroll position to match the saved view ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D60").Select()
